$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column header names (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B17").Value = "Bejucal De Ocampo"
$ws.Range("B18").Value = "Benemérito De Las Américas"
$ws.Range("B19").Value = "Comitán De Domínguez"
$ws.Range("B25").Value = "Mazapa De Madero"
$ws.Range("B29").Value = "Ocozocoautla De Espinosa"
$ws.Range("B31").Value = "Salto De Agua"
$ws.Range("A46").Value = "Ciudad De México"
$ws.Range("A53").Value = "Coahuila De Zaragoza"
$ws.Range("B58").Value = "San Juan Del Río"
$ws.Range("A61").Value = "Estado De México"
$ws.Range("B61").Value = "Almoloya De Juárez"
$ws.Range("B64").Value = "Ecatepec De Morelos"
$ws.Range("B69").Value = "Tlalnepantla De Baz"
$ws.Range("B80").Value = "Purísima Del Rincón"
$ws.Range("B83").Value = "San Francisco Del Rincón"
$ws.Range("B85").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B90").Value = "Acapulco De Juárez"
$ws.Range("B91").Value = "Ajuchitlán Del Progreso"
$ws.Range("B92").Value = "Atoyac De Álvarez"
$ws.Range("B93").Value = "Ayutla De Los Libres"
$ws.Range("B96").Value = "Coyuca De Catalán"
$ws.Range("B105").Value = "Técpan De Galeana"
$ws.Range("B106").Value = "Zihuatanejo De Azueta"
$ws.Range("B109").Value = "Pachuca De Soto"
$ws.Range("B110").Value = "Tula De Allende"
$ws.Range("B111").Value = "Tulancingo De Bravo"
$ws.Range("B114").Value = "Atotonilco El Alto"
$ws.Range("B115").Value = "Autlán De Navarro"
$ws.Range("B126").Value = "San Diego De Alejandría"
$ws.Range("B127").Value = "San Juan De Los Lagos"
$ws.Range("B129").Value = "Tamazula De Gordiano"
$ws.Range("B132").Value = "Tepatitlán De Morelos"
$ws.Range("B133").Value = "Tlajomulco De Zúñiga"
$ws.Range("B135").Value = "Unión De San Antonio"
$ws.Range("B136").Value = "Zacoalco De Torres"
$ws.Range("A140").Value = "Michoacán De Ocampo"
$ws.Range("B172").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B173").Value = "Guevea De Humboldt"
$ws.Range("B174").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B175").Value = "Huajuapan De León"
$ws.Range("B207").Value = "Santa María Del Río"
$ws.Range("B209").Value = "Villa De Arista"
$ws.Range("B210").Value = "Villa De Ramos"
$ws.Range("A227").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B235").Value = "Ignacio De La Llave"
$ws.Range("B238").Value = "Martínez De La Torre"
$ws.Range("B240").Value = "Mixtla De Altamirano"
$ws.Range("B243").Value = "Paso De Ovejas"
$ws.Range("B245").Value = "Sayula De Alemán"
$ws.Range("B247").Value = "Soledad De Doblado"
$ws.Range("B258").Value = "Trinidad García De La Cadena"
$ws.Range("B260").Value = "Villa De Cos"
$ws.Range("A263").Value = "Total"

# Remove trailing footnote rows (265-269), which also shrinks dimension to D263
$ws.Range("A265:A269").EntireRow.Delete()
